# Applies the "Updated cryptos list" data refresh to Sheet1 (rows 2-51).
# For each touched cell we force Text storage (NumberFormat "@") before writing
# the literal string, then ClearFormats() so the cell keeps default style 0 --
# this mirrors the source workbook, where every data cell is an unstyled
# inline string (prices like "538.36" must stay text, not become numeric).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "60.256.47"
Set-TextValue "E2" "  -2.41%  "

# Row 3
Set-TextValue "D3" "2.329.27"
Set-TextValue "E3" "  -6.03%  "

# Row 4
Set-TextValue "E4" "  +0.00%  "

# Row 5
Set-TextValue "D5" "538.36"
Set-TextValue "E5" "  -3.02%  "

# Row 6
Set-TextValue "D6" "135.93"
Set-TextValue "E6" "  -8.14%  "

# Row 7
Set-TextValue "E7" "  +0.02%  "

# Row 8
Set-TextValue "D8" "0.511"
Set-TextValue "E8" "  -15.11%  "

# Row 9
Set-TextValue "D9" "2.321.79"
Set-TextValue "E9" "  -6.29%  "

# Row 10
Set-TextValue "E10" "  -6.23%  "

# Row 11
Set-TextValue "E11" "  -0.44%  "

# Row 12
Set-TextValue "D12" "5.16"
Set-TextValue "E12" "  -6.23%  "

# Row 13
Set-TextValue "D13" "0.333"
Set-TextValue "E13" "  -6.87%  "

# Row 14
Set-TextValue "D14" "24.35"
Set-TextValue "E14" "  -8.36%  "

# Row 15
Set-TextValue "D15" "2.751.20"
Set-TextValue "E15" "  -5.98%  "

# Row 16
Set-TextValue "D16" "60.549.07"
Set-TextValue "E16" "  -1.82%  "

# Row 17
Set-TextValue "D17" "0.0000157"
Set-TextValue "E17" "  -6.90%  "

# Row 18
Set-TextValue "D18" "2.334.20"
Set-TextValue "E18" "  -5.88%  "

# Row 19
Set-TextValue "D19" "10.40"
Set-TextValue "E19" "  -7.59%  "

# Row 20
Set-TextValue "D20" "4.02"
Set-TextValue "E20" "  -5.12%  "

# Row 21
Set-TextValue "D21" "309.79"
Set-TextValue "E21" "  -4.07%  "

# Row 22
Set-TextValue "D22" "6.42"
Set-TextValue "E22" "  -11.35%  "

# Row 23
Set-TextValue "E23" "  -0.31%  "

# Row 24
Set-TextValue "E24" "  -3.13%  "

# Row 25
Set-TextValue "D25" "62.47"
Set-TextValue "E25" "  -3.17%  "

# Row 26
Set-TextValue "B26" "Binance-PegBSC-USD"
Set-TextValue "C26" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  +0.13%  "

# Row 27
Set-TextValue "B27" "Aptos"
Set-TextValue "C27" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D27" "7.94"
Set-TextValue "E27" "  +1.11%  "

# Row 28
Set-TextValue "D28" "2.449.32"
Set-TextValue "E28" "  -6.20%  "

# Row 29
Set-TextValue "D29" "0.0₃0874"
Set-TextValue "E29" "  -13.65%  "

# Row 30
Set-TextValue "B30" "InternetComputer(DFINITY)"
Set-TextValue "C30" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D30" "7.79"
Set-TextValue "E30" "  -7.44%  "

# Row 31
Set-TextValue "B31" "Bittensor"
Set-TextValue "C31" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D31" "496.99"
Set-TextValue "E31" "  -12.91%  "

# Row 32
Set-TextValue "D32" "1.36"
Set-TextValue "E32" "  -10.45%  "

# Row 33
Set-TextValue "D33" "0.142"
Set-TextValue "E33" "  -5.95%  "

# Row 34
Set-TextValue "D34" "1.77"
Set-TextValue "E34" "  -8.77%  "

# Row 35
Set-TextValue "D35" "1.50"
Set-TextValue "E35" "  -7.19%  "

# Row 36
Set-TextValue "D36" "0.998"
Set-TextValue "E36" "  -0.11%  "

# Row 37
Set-TextValue "B37" "PolygonEcosystemToken"
Set-TextValue "C37" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D37" "0.367"
Set-TextValue "E37" "  -4.52%  "

# Row 38
Set-TextValue "B38" "NEARProtocol"
Set-TextValue "C38" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D38" "4.44"
Set-TextValue "E38" "  -10.21%  "

# Row 39
Set-TextValue "D39" "18.09"
Set-TextValue "E39" "  -2.73%  "

# Row 40
Set-TextValue "D40" "5.17"
Set-TextValue "E40" "  -13.97%  "

# Row 41
Set-TextValue "E41" "  -1.56%  "

# Row 42
Set-TextValue "B42" "USDe"
Set-TextValue "C42" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D42" "1.00"
Set-TextValue "E42" "  +0.03%  "

# Row 43
Set-TextValue "B43" "Monero"
Set-TextValue "C43" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D43" "138.22"
Set-TextValue "E43" "  -4.24%  "

# Row 44
Set-TextValue "D44" "39.89"
Set-TextValue "E44" "  -1.91%  "

# Row 45
Set-TextValue "B45" "Aave"
Set-TextValue "C45" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D45" "135.54"
Set-TextValue "E45" "  -8.90%  "

# Row 46
Set-TextValue "B46" "Filecoin"
Set-TextValue "C46" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D46" "3.47"
Set-TextValue "E46" "  -4.90%  "

# Row 47
Set-TextValue "D47" "2.02"
Set-TextValue "E47" "  -17.72%  "

# Row 48
Set-TextValue "D48" "0.0503"
Set-TextValue "E48" "  -7.63%  "

# Row 49
Set-TextValue "D49" "19.10"
Set-TextValue "E49" "  -13.38%  "

# Row 50
Set-TextValue "D50" "0.560"
Set-TextValue "E50" "  -6.58%  "

# Row 51
Set-TextValue "D51" "0.0885"
Set-TextValue "E51" "  -6.35%  "
